$d = $word.ActiveDocument

# Replace all standalone "In;" paragraphs with "Wufake;"
$d.Content.Find.Execute("In;", $false, $true, $false, $false, $false, $true, 1, $false, "Wufake;", 2)

# Replace all "and out;" paragraphs with "wukhuphe;"
$d.Content.Find.Execute("and out;", $false, $true, $false, $false, $false, $true, 1, $false, "wukhuphe;", 2)
